$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 400
$ws.Range("I10").Value = 400
$ws.Range("K10").Value = 400
$ws.Range("M10").Value = -107

$ws.Range("H19").Value = 999.5
$ws.Range("I19").Value = 999.5
$ws.Range("K19").Value = 999.5
$ws.Range("M19").Value = -824.5

$ws.Range("H41").Value = 470.57144
$ws.Range("I41").Value = 429.76923
$ws.Range("K41").Value = 429.76923
$ws.Range("M41").Value = 10.23077000000001

$ws.Range("H70").Value = 3032.9167
$ws.Range("I70").Value = 2655
$ws.Range("K70").Value = 7965
$ws.Range("M70").Value = -7695

$ws.Range("H73").Value = 3032.9167
$ws.Range("I73").Value = 2655
$ws.Range("K73").Value = 7965
$ws.Range("M73").Value = -7029

$ws.Range("H86").Value = 1477.6
$ws.Range("I86").Value = 1664
$ws.Range("J86").Value = 1353.3334
$ws.Range("K86").Value = 1664
$ws.Range("L86").Value = 1353.3334
$ws.Range("M86").Value = -541
$ws.Range("N86").Value = -3599.3334

$ws.Range("H89").Value = 1477.6
$ws.Range("I89").Value = 1664
$ws.Range("J89").Value = 1353.3334
$ws.Range("K89").Value = 8320
$ws.Range("L89").Value = 6766.666999999999
$ws.Range("M89").Value = -2704
$ws.Range("N89").Value = -17998.667

$ws.Range("H100").Value = 3031.6667
$ws.Range("I100").Value = 2100
$ws.Range("J100").Value = 3497.5
$ws.Range("K100").Value = 2100
$ws.Range("L100").Value = 3497.5
$ws.Range("M100").Value = -1559
$ws.Range("N100").Value = -4579.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 9999
$ws.Range("I6").Value = 9999
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 9999
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -9826
$ws.Range("N6").ClearContents()

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H86").Value = 2494.9167
$ws.Range("I86").Value = 2093.9
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 2093.9
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -970.9000000000001
$ws.Range("N86").Value = -6746

$ws.Range("H89").Value = 2494.9167
$ws.Range("I89").Value = 2093.9
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 10469.5
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -4853.5
$ws.Range("N89").Value = -33732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 283.2353
$ws.Range("I7").Value = 269.8125
$ws.Range("J7").Value = 498
$ws.Range("K7").Value = 269.8125
$ws.Range("L7").Value = 498
$ws.Range("M7").Value = -156.8125
$ws.Range("N7").Value = -724

$ws.Range("H12").Value = 2102
$ws.Range("I12").Value = 381.6
$ws.Range("J12").Value = 6403
$ws.Range("K12").Value = 381.6
$ws.Range("L12").Value = 6403
$ws.Range("M12").Value = -211.6
$ws.Range("N12").Value = -6743

$ws.Range("H31").Value = 3398.9
$ws.Range("I31").Value = 1158.1428
$ws.Range("K31").Value = 1158.1428
$ws.Range("M31").Value = -863.1428000000001

$ws.Range("H34").Value = 3398.9
$ws.Range("I34").Value = 1158.1428
$ws.Range("K34").Value = 1158.1428
$ws.Range("M34").Value = -956.1428000000001

$ws.Range("H105").Value = 1563.4
$ws.Range("I105").Value = 1563.4
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1563.4
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 183.5999999999999
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 141.66667
$ws.Range("I17").Value = 141.66667
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 425.00001
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -256.00001
$ws.Range("N17").ClearContents()

$ws.Range("H64").Value = 1418.5
$ws.Range("I64").Value = 982.2
$ws.Range("K64").Value = 2946.6
$ws.Range("M64").Value = -2676.6

$ws.Range("H67").Value = 1418.5
$ws.Range("I67").Value = 982.2
$ws.Range("K67").Value = 2946.6
$ws.Range("M67").Value = -2010.6

$ws.Range("H98").Value = 255
$ws.Range("I98").Value = 255
$ws.Range("K98").Value = 765
$ws.Range("M98").Value = 733

$ws.Range("H114").Value = 1125.3334
$ws.Range("I114").Value = 980.75
$ws.Range("J114").Value = 1197.625
$ws.Range("K114").Value = 2942.25
$ws.Range("L114").Value = 3592.875
$ws.Range("M114").Value = 311.75
$ws.Range("N114").Value = -10100.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 539.8
$ws.Range("I2").Value = 706.8570999999999
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 706.8570999999999
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -593.8570999999999
$ws.Range("N2").Value = -376

$ws.Range("H70").Value = 83338200
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 83338200
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 83338200
$ws.Range("N70").Value = -83338740
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 83338200
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 83338200
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 83338200
$ws.Range("N73").Value = -83340072
$ws.Range("M73").ClearContents()

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1301.5
$ws.Range("I16").Value = 1326
$ws.Range("J16").Value = 1179
$ws.Range("K16").Value = 1326
$ws.Range("L16").Value = 1179
$ws.Range("M16").Value = -1156
$ws.Range("N16").Value = -1519

$ws.Range("H136").Value = 2422.1538
$ws.Range("I136").Value = 2284.8572
$ws.Range("J136").Value = 2582.3333
$ws.Range("K136").Value = 6854.571599999999
$ws.Range("L136").Value = 7746.999899999999
$ws.Range("M136").Value = -4304.571599999999
$ws.Range("N136").Value = -12846.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 8335400
$ws.Range("I23").Value = 466.66666
$ws.Range("J23").Value = 16670333
$ws.Range("K23").Value = 466.66666
$ws.Range("L23").Value = 16670333
$ws.Range("M23").Value = -237.66666
$ws.Range("N23").Value = -16670791

$ws.Range("H136").Value = 809.4545000000001
$ws.Range("I136").Value = 809.4545000000001
$ws.Range("K136").Value = 2428.3635
$ws.Range("M136").Value = 121.6364999999996
